$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Header field updates ---
$ws.Cells.Item(11,5).Value = 1700858
$ws.Cells.Item(13,3).Value = 7
$ws.Cells.Item(13,6).Value = 9

# --- Restructure data table: insert one row at 63 so the footer (rows 67-68) shifts to 68-69 ---
$ws.Rows(63).Insert()

# Fix up borders: row 63 (new) should get the "last row" bottom-border style currently sitting on row 62,
# and row 62 should go back to being a normal interior row (style currently on row 61).
$ws.Range("B62:J62").Copy()
$ws.Range("B63:J63").PasteSpecial(-4122)
$ws.Range("B61:J61").Copy()
$ws.Range("B62:J62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Write the reordered / updated employee data (grouped by employee, periods descending) ---
$ws.Cells.Item(16,2).Value = 'CC'
$ws.Cells.Item(16,3).Value = '73196440'
$ws.Cells.Item(16,4).Value = 'ARMANDO DE JESUS BLANCO VEGA'
$ws.Cells.Item(16,5).Value = '2208'
$ws.Cells.Item(16,6).Value = 26650
$ws.Cells.Item(16,7).Value = 908526
$ws.Cells.Item(17,2).Value = 'CC'
$ws.Cells.Item(17,3).Value = '73196440'
$ws.Cells.Item(17,4).Value = 'ARMANDO DE JESUS BLANCO VEGA'
$ws.Cells.Item(17,5).Value = '2207'
$ws.Cells.Item(17,6).Value = 36341
$ws.Cells.Item(17,7).Value = 908526
$ws.Cells.Item(18,2).Value = 'CC'
$ws.Cells.Item(18,3).Value = '73196440'
$ws.Cells.Item(18,4).Value = 'ARMANDO DE JESUS BLANCO VEGA'
$ws.Cells.Item(18,5).Value = '2206'
$ws.Cells.Item(18,6).Value = 36341
$ws.Cells.Item(18,7).Value = 908526
$ws.Cells.Item(19,2).Value = 'CC'
$ws.Cells.Item(19,3).Value = '73196440'
$ws.Cells.Item(19,4).Value = 'ARMANDO DE JESUS BLANCO VEGA'
$ws.Cells.Item(19,5).Value = '2205'
$ws.Cells.Item(19,6).Value = 36341
$ws.Cells.Item(19,7).Value = 908526
$ws.Cells.Item(20,2).Value = 'CC'
$ws.Cells.Item(20,3).Value = '73196440'
$ws.Cells.Item(20,4).Value = 'ARMANDO DE JESUS BLANCO VEGA'
$ws.Cells.Item(20,5).Value = '2204'
$ws.Cells.Item(20,6).Value = 36341
$ws.Cells.Item(20,7).Value = 908526
$ws.Cells.Item(21,2).Value = 'CC'
$ws.Cells.Item(21,3).Value = '73196440'
$ws.Cells.Item(21,4).Value = 'ARMANDO DE JESUS BLANCO VEGA'
$ws.Cells.Item(21,5).Value = '2203'
$ws.Cells.Item(21,6).Value = 36341
$ws.Cells.Item(21,7).Value = 908526
$ws.Cells.Item(22,2).Value = 'CC'
$ws.Cells.Item(22,3).Value = '73196440'
$ws.Cells.Item(22,4).Value = 'ARMANDO DE JESUS BLANCO VEGA'
$ws.Cells.Item(22,5).Value = '2202'
$ws.Cells.Item(22,6).Value = 36341
$ws.Cells.Item(22,7).Value = 908526
$ws.Cells.Item(23,2).Value = 'CC'
$ws.Cells.Item(23,3).Value = '73196440'
$ws.Cells.Item(23,4).Value = 'ARMANDO DE JESUS BLANCO VEGA'
$ws.Cells.Item(23,5).Value = '2201'
$ws.Cells.Item(23,6).Value = 36341
$ws.Cells.Item(23,7).Value = 908526
$ws.Cells.Item(24,2).Value = 'CC'
$ws.Cells.Item(24,3).Value = '1128054471'
$ws.Cells.Item(24,4).Value = 'MANUEL ALBERTO OSORIO SALGADO'
$ws.Cells.Item(24,5).Value = '2208'
$ws.Cells.Item(24,6).Value = 26650
$ws.Cells.Item(24,7).Value = 908526
$ws.Cells.Item(25,2).Value = 'CC'
$ws.Cells.Item(25,3).Value = '1128054471'
$ws.Cells.Item(25,4).Value = 'MANUEL ALBERTO OSORIO SALGADO'
$ws.Cells.Item(25,5).Value = '2207'
$ws.Cells.Item(25,6).Value = 36341
$ws.Cells.Item(25,7).Value = 908526
$ws.Cells.Item(26,2).Value = 'CC'
$ws.Cells.Item(26,3).Value = '1128054471'
$ws.Cells.Item(26,4).Value = 'MANUEL ALBERTO OSORIO SALGADO'
$ws.Cells.Item(26,5).Value = '2206'
$ws.Cells.Item(26,6).Value = 36341
$ws.Cells.Item(26,7).Value = 908526
$ws.Cells.Item(27,2).Value = 'CC'
$ws.Cells.Item(27,3).Value = '1128054471'
$ws.Cells.Item(27,4).Value = 'MANUEL ALBERTO OSORIO SALGADO'
$ws.Cells.Item(27,5).Value = '2205'
$ws.Cells.Item(27,6).Value = 36341
$ws.Cells.Item(27,7).Value = 908526
$ws.Cells.Item(28,2).Value = 'CC'
$ws.Cells.Item(28,3).Value = '1128054471'
$ws.Cells.Item(28,4).Value = 'MANUEL ALBERTO OSORIO SALGADO'
$ws.Cells.Item(28,5).Value = '2204'
$ws.Cells.Item(28,6).Value = 36341
$ws.Cells.Item(28,7).Value = 908526
$ws.Cells.Item(29,2).Value = 'CC'
$ws.Cells.Item(29,3).Value = '1128054471'
$ws.Cells.Item(29,4).Value = 'MANUEL ALBERTO OSORIO SALGADO'
$ws.Cells.Item(29,5).Value = '2203'
$ws.Cells.Item(29,6).Value = 36341
$ws.Cells.Item(29,7).Value = 908526
$ws.Cells.Item(30,2).Value = 'CC'
$ws.Cells.Item(30,3).Value = '1128054471'
$ws.Cells.Item(30,4).Value = 'MANUEL ALBERTO OSORIO SALGADO'
$ws.Cells.Item(30,5).Value = '2202'
$ws.Cells.Item(30,6).Value = 36341
$ws.Cells.Item(30,7).Value = 908526
$ws.Cells.Item(31,2).Value = 'CC'
$ws.Cells.Item(31,3).Value = '1128054471'
$ws.Cells.Item(31,4).Value = 'MANUEL ALBERTO OSORIO SALGADO'
$ws.Cells.Item(31,5).Value = '2201'
$ws.Cells.Item(31,6).Value = 36341
$ws.Cells.Item(31,7).Value = 908526
$ws.Cells.Item(32,2).Value = 'CC'
$ws.Cells.Item(32,3).Value = '1007255013'
$ws.Cells.Item(32,4).Value = 'NAUDITH PALOMINO ARIZA'
$ws.Cells.Item(32,5).Value = '2208'
$ws.Cells.Item(32,6).Value = 26650
$ws.Cells.Item(32,7).Value = 1000000
$ws.Cells.Item(33,2).Value = 'CC'
$ws.Cells.Item(33,3).Value = '1007255013'
$ws.Cells.Item(33,4).Value = 'NAUDITH PALOMINO ARIZA'
$ws.Cells.Item(33,5).Value = '2207'
$ws.Cells.Item(33,6).Value = 40000
$ws.Cells.Item(33,7).Value = 1000000
$ws.Cells.Item(34,2).Value = 'CC'
$ws.Cells.Item(34,3).Value = '1007255013'
$ws.Cells.Item(34,4).Value = 'NAUDITH PALOMINO ARIZA'
$ws.Cells.Item(34,5).Value = '2206'
$ws.Cells.Item(34,6).Value = 40000
$ws.Cells.Item(34,7).Value = 1000000
$ws.Cells.Item(35,2).Value = 'CC'
$ws.Cells.Item(35,3).Value = '1007255013'
$ws.Cells.Item(35,4).Value = 'NAUDITH PALOMINO ARIZA'
$ws.Cells.Item(35,5).Value = '2205'
$ws.Cells.Item(35,6).Value = 40000
$ws.Cells.Item(35,7).Value = 1000000
$ws.Cells.Item(36,2).Value = 'CC'
$ws.Cells.Item(36,3).Value = '1007255013'
$ws.Cells.Item(36,4).Value = 'NAUDITH PALOMINO ARIZA'
$ws.Cells.Item(36,5).Value = '2204'
$ws.Cells.Item(36,6).Value = 40000
$ws.Cells.Item(36,7).Value = 1000000
$ws.Cells.Item(37,2).Value = 'CC'
$ws.Cells.Item(37,3).Value = '1007255013'
$ws.Cells.Item(37,4).Value = 'NAUDITH PALOMINO ARIZA'
$ws.Cells.Item(37,5).Value = '2202'
$ws.Cells.Item(37,6).Value = 36341
$ws.Cells.Item(37,7).Value = 1000000
$ws.Cells.Item(38,2).Value = 'CC'
$ws.Cells.Item(38,3).Value = '1007255013'
$ws.Cells.Item(38,4).Value = 'NAUDITH PALOMINO ARIZA'
$ws.Cells.Item(38,5).Value = '2201'
$ws.Cells.Item(38,6).Value = 36341
$ws.Cells.Item(38,7).Value = 1000000
$ws.Cells.Item(39,2).Value = 'CC'
$ws.Cells.Item(39,3).Value = '1047498248'
$ws.Cells.Item(39,4).Value = 'JORGE LUIS MORELOS PEREZ'
$ws.Cells.Item(39,5).Value = '2208'
$ws.Cells.Item(39,6).Value = 26650
$ws.Cells.Item(39,7).Value = 871000
$ws.Cells.Item(40,2).Value = 'CC'
$ws.Cells.Item(40,3).Value = '1047498248'
$ws.Cells.Item(40,4).Value = 'JORGE LUIS MORELOS PEREZ'
$ws.Cells.Item(40,5).Value = '2207'
$ws.Cells.Item(40,6).Value = 36341
$ws.Cells.Item(40,7).Value = 871000
$ws.Cells.Item(41,2).Value = 'CC'
$ws.Cells.Item(41,3).Value = '1047498248'
$ws.Cells.Item(41,4).Value = 'JORGE LUIS MORELOS PEREZ'
$ws.Cells.Item(41,5).Value = '2206'
$ws.Cells.Item(41,6).Value = 36341
$ws.Cells.Item(41,7).Value = 871000
$ws.Cells.Item(42,2).Value = 'CC'
$ws.Cells.Item(42,3).Value = '1047498248'
$ws.Cells.Item(42,4).Value = 'JORGE LUIS MORELOS PEREZ'
$ws.Cells.Item(42,5).Value = '2205'
$ws.Cells.Item(42,6).Value = 36341
$ws.Cells.Item(42,7).Value = 871000
$ws.Cells.Item(43,2).Value = 'CC'
$ws.Cells.Item(43,3).Value = '1047498248'
$ws.Cells.Item(43,4).Value = 'JORGE LUIS MORELOS PEREZ'
$ws.Cells.Item(43,5).Value = '2204'
$ws.Cells.Item(43,6).Value = 36341
$ws.Cells.Item(43,7).Value = 871000
$ws.Cells.Item(44,2).Value = 'CC'
$ws.Cells.Item(44,3).Value = '1047498248'
$ws.Cells.Item(44,4).Value = 'JORGE LUIS MORELOS PEREZ'
$ws.Cells.Item(44,5).Value = '2203'
$ws.Cells.Item(44,6).Value = 36341
$ws.Cells.Item(44,7).Value = 871000
$ws.Cells.Item(45,2).Value = 'CC'
$ws.Cells.Item(45,3).Value = '1047498248'
$ws.Cells.Item(45,4).Value = 'JORGE LUIS MORELOS PEREZ'
$ws.Cells.Item(45,5).Value = '2202'
$ws.Cells.Item(45,6).Value = 36341
$ws.Cells.Item(45,7).Value = 871000
$ws.Cells.Item(46,2).Value = 'CC'
$ws.Cells.Item(46,3).Value = '1047498248'
$ws.Cells.Item(46,4).Value = 'JORGE LUIS MORELOS PEREZ'
$ws.Cells.Item(46,5).Value = '2201'
$ws.Cells.Item(46,6).Value = 36341
$ws.Cells.Item(46,7).Value = 871000
$ws.Cells.Item(47,2).Value = 'CC'
$ws.Cells.Item(47,3).Value = '22800429'
$ws.Cells.Item(47,4).Value = 'ELIDA OSPINO VASQUEZ'
$ws.Cells.Item(47,5).Value = '2208'
$ws.Cells.Item(47,6).Value = 26650
$ws.Cells.Item(47,7).Value = 908526
$ws.Cells.Item(48,2).Value = 'CC'
$ws.Cells.Item(48,3).Value = '22800429'
$ws.Cells.Item(48,4).Value = 'ELIDA OSPINO VASQUEZ'
$ws.Cells.Item(48,5).Value = '2207'
$ws.Cells.Item(48,6).Value = 36341
$ws.Cells.Item(48,7).Value = 908526
$ws.Cells.Item(49,2).Value = 'CC'
$ws.Cells.Item(49,3).Value = '22800429'
$ws.Cells.Item(49,4).Value = 'ELIDA OSPINO VASQUEZ'
$ws.Cells.Item(49,5).Value = '2206'
$ws.Cells.Item(49,6).Value = 36341
$ws.Cells.Item(49,7).Value = 908526
$ws.Cells.Item(50,2).Value = 'CC'
$ws.Cells.Item(50,3).Value = '22800429'
$ws.Cells.Item(50,4).Value = 'ELIDA OSPINO VASQUEZ'
$ws.Cells.Item(50,5).Value = '2205'
$ws.Cells.Item(50,6).Value = 36341
$ws.Cells.Item(50,7).Value = 908526
$ws.Cells.Item(51,2).Value = 'CC'
$ws.Cells.Item(51,3).Value = '22800429'
$ws.Cells.Item(51,4).Value = 'ELIDA OSPINO VASQUEZ'
$ws.Cells.Item(51,5).Value = '2204'
$ws.Cells.Item(51,6).Value = 36341
$ws.Cells.Item(51,7).Value = 908526
$ws.Cells.Item(52,2).Value = 'CC'
$ws.Cells.Item(52,3).Value = '22800429'
$ws.Cells.Item(52,4).Value = 'ELIDA OSPINO VASQUEZ'
$ws.Cells.Item(52,5).Value = '2203'
$ws.Cells.Item(52,6).Value = 36341
$ws.Cells.Item(52,7).Value = 908526
$ws.Cells.Item(53,2).Value = 'CC'
$ws.Cells.Item(53,3).Value = '22800429'
$ws.Cells.Item(53,4).Value = 'ELIDA OSPINO VASQUEZ'
$ws.Cells.Item(53,5).Value = '2202'
$ws.Cells.Item(53,6).Value = 36341
$ws.Cells.Item(53,7).Value = 908526
$ws.Cells.Item(54,2).Value = 'CC'
$ws.Cells.Item(54,3).Value = '22800429'
$ws.Cells.Item(54,4).Value = 'ELIDA OSPINO VASQUEZ'
$ws.Cells.Item(54,5).Value = '2201'
$ws.Cells.Item(54,6).Value = 36341
$ws.Cells.Item(54,7).Value = 908526
$ws.Cells.Item(55,2).Value = 'CC'
$ws.Cells.Item(55,3).Value = '3800855'
$ws.Cells.Item(55,4).Value = 'WILLIAM DIAZ MARTINEZ'
$ws.Cells.Item(55,5).Value = '2208'
$ws.Cells.Item(55,6).Value = 26650
$ws.Cells.Item(55,7).Value = 908526
$ws.Cells.Item(56,2).Value = 'CC'
$ws.Cells.Item(56,3).Value = '3800855'
$ws.Cells.Item(56,4).Value = 'WILLIAM DIAZ MARTINEZ'
$ws.Cells.Item(56,5).Value = '2207'
$ws.Cells.Item(56,6).Value = 36341
$ws.Cells.Item(56,7).Value = 908526
$ws.Cells.Item(57,2).Value = 'CC'
$ws.Cells.Item(57,3).Value = '3800855'
$ws.Cells.Item(57,4).Value = 'WILLIAM DIAZ MARTINEZ'
$ws.Cells.Item(57,5).Value = '2206'
$ws.Cells.Item(57,6).Value = 36341
$ws.Cells.Item(57,7).Value = 908526
$ws.Cells.Item(58,2).Value = 'CC'
$ws.Cells.Item(58,3).Value = '3800855'
$ws.Cells.Item(58,4).Value = 'WILLIAM DIAZ MARTINEZ'
$ws.Cells.Item(58,5).Value = '2205'
$ws.Cells.Item(58,6).Value = 36341
$ws.Cells.Item(58,7).Value = 908526
$ws.Cells.Item(59,2).Value = 'CC'
$ws.Cells.Item(59,3).Value = '3800855'
$ws.Cells.Item(59,4).Value = 'WILLIAM DIAZ MARTINEZ'
$ws.Cells.Item(59,5).Value = '2204'
$ws.Cells.Item(59,6).Value = 36341
$ws.Cells.Item(59,7).Value = 908526
$ws.Cells.Item(60,2).Value = 'CC'
$ws.Cells.Item(60,3).Value = '3800855'
$ws.Cells.Item(60,4).Value = 'WILLIAM DIAZ MARTINEZ'
$ws.Cells.Item(60,5).Value = '2203'
$ws.Cells.Item(60,6).Value = 36341
$ws.Cells.Item(60,7).Value = 908526
$ws.Cells.Item(61,2).Value = 'CC'
$ws.Cells.Item(61,3).Value = '3800855'
$ws.Cells.Item(61,4).Value = 'WILLIAM DIAZ MARTINEZ'
$ws.Cells.Item(61,5).Value = '2202'
$ws.Cells.Item(61,6).Value = 36341
$ws.Cells.Item(61,7).Value = 908526
$ws.Cells.Item(62,2).Value = 'CC'
$ws.Cells.Item(62,3).Value = '3800855'
$ws.Cells.Item(62,4).Value = 'WILLIAM DIAZ MARTINEZ'
$ws.Cells.Item(62,5).Value = '2201'
$ws.Cells.Item(62,6).Value = 36341
$ws.Cells.Item(62,7).Value = 908526
$ws.Cells.Item(63,2).Value = 'CC'
$ws.Cells.Item(63,3).Value = '1143397931'
$ws.Cells.Item(63,4).Value = 'ELKIN DE JESUS VILLA NUÑEZ'
$ws.Cells.Item(63,5).Value = '2106'
$ws.Cells.Item(63,6).Value = 36341
$ws.Cells.Item(63,7).Value = 908526

# --- Footer already shifted by the row insert above: row68 = signature line, row69 = name/firma labels ---
